# Update the "lifts" worksheet state table:
#  - change row 5 SITE -> WAREHOUSE
#  - append a new row 6 with a SHORE entry

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lifts")

# Update existing cell B5: SITE -> WAREHOUSE
$ws.Range("B5").Value = "WAREHOUSE"

# Add the new row 6 values (unstyled, same as row 5)
$ws.Cells.Item(6, 1).Value = 3
$ws.Cells.Item(6, 2).Value = "SHORE"
$ws.Cells.Item(6, 3).Value = "S"
$ws.Cells.Item(6, 4).Value = "O"
$ws.Cells.Item(6, 5).Value = "N"

$ws.Range("A6:E6").Style = "Normal"
